$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update shared string "polygon" -> "bnb" (cell E2)
$ws.Range("E2").Value = "bnb"

# 2. Update long instructional text in H4
$ws.Range("H4").Value = "< - initialiser la première ligne manuellement pour commencer immédiatement =  ////200000 //// 0,00000000001 //// date du jour //// heure"

# 3. Update selection (active cell) to G34
$ws.Range("G34").Select()

# 4. Update column widths.
# The runtime persists <col width="..."/> as (round(ColumnWidth*6)+5)/6 -
# i.e. ColumnWidth is quantized to 1/6-character pixel steps (same rounding
# Excel itself applies at 96dpi/Calibri 11). We pick, for each column, the
# ColumnWidth whose quantized result lands closest to the diff's target
# stored width.
$ws.Columns.Item(2).ColumnWidth = 18.1666666666667
$ws.Columns.Item(3).ColumnWidth = 20.1666666666667
$ws.Columns.Item(4).ColumnWidth = 23.1666666666667
$ws.Columns.Item(5).ColumnWidth = 24
$ws.Columns.Item(6).ColumnWidth = 14.3333333333333
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 16.6666666666667
$ws.Columns.Item(9).ColumnWidth = 9.16666666666667
$ws.Columns.Item(11).ColumnWidth = 9.66666666666667
$ws.Columns.Item(12).ColumnWidth = 8
$ws.Columns.Item(13).ColumnWidth = 7.66666666666667
$ws.Columns.Item(14).ColumnWidth = 13.3333333333333

# 5. Update numeric values (A4, and dependent formula cells I2, J2, K2 recompute automatically)
$ws.Range("A4").Value = 200000
